$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $old with $new inside $shape's
# TextRange, splicing only the matched characters so that surrounding runs
# (and their formatting, e.g. red-colored substrings) are left untouched.
# NOTE: call helpers with POSITIONAL arguments only - named arguments
# (-foo bar) do not bind correctly with this host's PowerShell engine.
# ---------------------------------------------------------------------------
function Replace-InShapeText {
    param($shape, [string]$old, [string]$new)
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($old)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $old.Length)
        $sub.Text = $new
        return $true
    }
    return $false
}

# Find the first shape on a slide whose text contains $needle.
function Find-ShapeWithText {
    param($slide, [string]$needle)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text.Contains($needle)) {
                    return $shp
                }
            }
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Delete slide 11 ("通報機能"). The old slide 12 ("デモプレイ") shifts up
#    and becomes the new slide 11 - its content is unchanged.
# ---------------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$reportShape = Find-ShapeWithText $slide11 "通報機能"
if ($reportShape -ne $null) {
    $slide11.Delete()
}

# ---------------------------------------------------------------------------
# 2) Slide (was 3) "SKET" slide: "Studentask to & evaluate Teacher"
#    -> "Studentask & evaluate Teacher"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp = Find-ShapeWithText $s3 " to & "
if ($shp -ne $null) {
    Replace-InShapeText $shp " to & " " & " | Out-Null
}

# ---------------------------------------------------------------------------
# 3) Slide (was 5) "本システムの特長": remove the standalone "通報機能"
#    paragraph, keeping the "グループ機能" paragraph right before it.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp = Find-ShapeWithText $s5 "通報"
if ($shp -ne $null) {
    $tr = $shp.TextFrame.TextRange
    $paraCount = $tr.Paragraphs().Count
    for ($i = $paraCount; $i -ge 1; $i--) {
        $para = $tr.Paragraphs($i, 1)
        if ($para.Text.Contains("通報")) {
            $para.Delete()
        }
    }
}

# ---------------------------------------------------------------------------
# 4) Slide (was 9): "（○○）を見ることができる。" -> "（解答数、評価値など）を見ることができる。"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp = Find-ShapeWithText $s9 "（○○）を見ることができる。"
if ($shp -ne $null) {
    Replace-InShapeText $shp "（○○）を見ることができる。" "（解答数、評価値など）を見ることができる。" | Out-Null
}

# ---------------------------------------------------------------------------
# 5) Slide (was 10): "内のみに向けた質問なども可能。" -> "内のみに向けた質問が可能。"
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp = Find-ShapeWithText $s10 "内のみに向けた質問なども可能。"
if ($shp -ne $null) {
    Replace-InShapeText $shp "内のみに向けた質問なども可能。" "内のみに向けた質問が可能。" | Out-Null
}

# ---------------------------------------------------------------------------
# 6) Date placeholder text cache bump: 2018/7/21 -> 2018/7/22
#    (slide master, every slide layout, and the notes master)
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
foreach ($shp in $master.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "2018/7/21") {
        $shp.TextFrame.TextRange.Text = "2018/7/22"
    }
}
foreach ($layout in $master.CustomLayouts) {
    foreach ($shp in $layout.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "2018/7/21") {
            $shp.TextFrame.TextRange.Text = "2018/7/22"
        }
    }
}

$notesMaster = $p.NotesMaster
$nmDate = $notesMaster.HeadersFooters.DateAndTime
$nmDate.Text = "2018/7/22"
